$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436

$ws.Range("D32").Value = 8423
$ws.Range("D33").Value = 9386
$ws.Range("D34").Value = 10671
$ws.Range("D35").Value = 13277
$ws.Range("D36").Value = 13758
$ws.Range("D37").Value = 14462
